$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction data for rows 2-7 (columns A-L)
$data = @(
    @(8,  1, "Cash To Mustafa", 1000, "MZN", 1,  1, 1000, "MZN", 2, "Clearance",       "money transfer to xyz"),
    @(9,  1, "Cash To Mustafa", 950,  "MZN", 1,  1, 950,  "MZN", 2, "Clearance",       "money transfer to xyz"),
    @(10, 1, "Cash To Mustafa", 999,  "MZN", 1,  1, 999,  "MZN", 2, "Clearance",       "money transfer to xyz"),
    @(11, 2, "Clearance",       1,    "MZN", 1,  1, 1,    "MZN", 1, "Cash To Mustafa", "money transfer to xyz"),
    @(12, 2, "Clearance",       888,  "MZN", 5,  1, 4440, "EUR", 1, "Cash To Mustafa", "money transfer to xyz"),
    @(13, 1, "Cash To Mustafa", 55,   "USD", 11, 1, 605,  "AED", 2, "Clearance",       "money transfer to xyz")
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
